$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-7) before writing new data (rows 2-10)
$ws.Range("A2:T10").ClearContents()

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Dll1"
$ws.Range("C2").Value = "Notch3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.190862666666667
$ws.Range("H2").Value = 15.572588
$ws.Range("I2").Value = 0.8740249884703439
$ws.Range("J2").Value = 0.874024988470344
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.847498666666667
$ws.Range("N2").Value = 14.542496
$ws.Range("O2").Value = 0.03400671694637637
$ws.Range("P2").Value = 0.03400671694637637
$ws.Range("Q2").Value = 25.16269985551644
$ws.Range("R2").Value = 226.464298699648
$ws.Range("S2").Value = 0.02972272038697085
$ws.Range("T2").Value = 0.02972272038697086

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Dll1"
$ws.Range("C3").Value = "Notch3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.190862666666667
$ws.Range("H3").Value = 15.572588
$ws.Range("I3").Value = 0.8740249884703439
$ws.Range("J3").Value = 0.874024988470344
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.627093333333334
$ws.Range("N3").Value = 4.88128
$ws.Range("O3").Value = 0.01141456784970118
$ws.Range("P3").Value = 0.01141456784970118
$ws.Range("Q3").Value = 8.446018039182224
$ws.Range("R3").Value = 76.01416235264001
$ws.Range("S3").Value = 0.009976617533229034
$ws.Range("T3").Value = 0.009976617533229036

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Dll1"
$ws.Range("C4").Value = "Notch3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.190862666666667
$ws.Range("H4").Value = 15.572588
$ws.Range("I4").Value = 0.8740249884703439
$ws.Range("J4").Value = 0.874024988470344
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 136.0707373333333
$ws.Range("N4").Value = 408.212212
$ws.Range("O4").Value = 0.9545787152039225
$ws.Range("P4").Value = 0.9545787152039225
$ws.Range("Q4").Value = 706.3245104494063
$ws.Range("R4").Value = 6356.920594044657
$ws.Range("S4").Value = 0.8343256505501441
$ws.Range("T4").Value = 0.8343256505501442

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Dll1"
$ws.Range("C5").Value = "Notch3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.100996
$ws.Range("H5").Value = 0.302988
$ws.Range("I5").Value = 0.0170054639091879
$ws.Range("J5").Value = 0.0170054639091879
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.847498666666667
$ws.Range("N5").Value = 14.542496
$ws.Range("O5").Value = 0.03400671694637637
$ws.Range("P5").Value = 0.03400671694637637
$ws.Range("Q5").Value = 0.4895779753386666
$ws.Range("R5").Value = 4.406201778048
$ws.Range("S5").Value = 0.0005782999977015718
$ws.Range("T5").Value = 0.0005782999977015718

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Dll1"
$ws.Range("C6").Value = "Notch3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.100996
$ws.Range("H6").Value = 0.302988
$ws.Range("I6").Value = 0.0170054639091879
$ws.Range("J6").Value = 0.0170054639091879
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.627093333333334
$ws.Range("N6").Value = 4.88128
$ws.Range("O6").Value = 0.01141456784970118
$ws.Range("P6").Value = 0.01141456784970118
$ws.Range("Q6").Value = 0.1643299182933333
$ws.Range("R6").Value = 1.47896926464
$ws.Range("S6").Value = 0.0001941100216070699
$ws.Range("T6").Value = 0.0001941100216070699

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Dll1"
$ws.Range("C7").Value = "Notch3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.100996
$ws.Range("H7").Value = 0.302988
$ws.Range("I7").Value = 0.0170054639091879
$ws.Range("J7").Value = 0.0170054639091879
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 136.0707373333333
$ws.Range("N7").Value = 408.212212
$ws.Range("O7").Value = 0.9545787152039225
$ws.Range("P7").Value = 0.9545787152039225
$ws.Range("Q7").Value = 13.74260018771733
$ws.Range("R7").Value = 123.683401689456
$ws.Range("S7").Value = 0.01623305388987925
$ws.Range("T7").Value = 0.01623305388987925

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Dll1"
$ws.Range("C8").Value = "Notch3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6471736666666666
$ws.Range("H8").Value = 1.941521
$ws.Range("I8").Value = 0.1089695476204681
$ws.Range("J8").Value = 0.1089695476204681
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.847498666666667
$ws.Range("N8").Value = 14.542496
$ws.Range("O8").Value = 0.03400671694637637
$ws.Range("P8").Value = 0.03400671694637637
$ws.Range("Q8").Value = 3.137173486268444
$ws.Range("R8").Value = 28.234561376416
$ws.Range("S8").Value = 0.00370569656170394
$ws.Range("T8").Value = 0.00370569656170394

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Dll1"
$ws.Range("C9").Value = "Notch3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6471736666666666
$ws.Range("H9").Value = 1.941521
$ws.Range("I9").Value = 0.1089695476204681
$ws.Range("J9").Value = 0.1089695476204681
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.627093333333334
$ws.Range("N9").Value = 4.88128
$ws.Range("O9").Value = 0.01141456784970118
$ws.Range("P9").Value = 0.01141456784970118
$ws.Range("Q9").Value = 1.053011958542222
$ws.Range("R9").Value = 9.47710762688
$ws.Range("S9").Value = 0.001243840294865077
$ws.Range("T9").Value = 0.001243840294865077

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Dll1"
$ws.Range("C10").Value = "Notch3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.6471736666666666
$ws.Range("H10").Value = 1.941521
$ws.Range("I10").Value = 0.1089695476204681
$ws.Range("J10").Value = 0.1089695476204681
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 136.0707373333333
$ws.Range("N10").Value = 408.212212
$ws.Range("O10").Value = 0.9545787152039225
$ws.Range("P10").Value = 0.9545787152039225
$ws.Range("Q10").Value = 88.06139800605023
$ws.Range("R10").Value = 792.552582054452
$ws.Range("S10").Value = 0.1040200107638991
$ws.Range("T10").Value = 0.1040200107638991
